# "Updated with Few more examples"
# Adds a new "RiskFactor-DirtectResponse" example to the Examples sheet:
#  - Inserts two new columns (RequestContent / ResponseContent) before the
#    existing ResponseByFields/Csvson columns.
#  - Appends two new trailing columns (StoreResponseVariables /
#    EvaluateFunctionVariables).
#  - Appends a new data row (row 6) describing the RiskFactor example.
#  - Adds hyperlinks for the new row's URL and for the pre-existing
#    bgates URL example (row 4) that didn't have one yet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Examples")

# --- 1. Insert two new columns at I:J (RequestContent / ResponseContent) ---
# This shifts the existing ResponseByFields (I) / Csvson (J) columns to K/L.
$ws.Range("I1:J1").EntireColumn.Insert()

# Header text for the two newly inserted columns.
$ws.Range("I1").Value = "RequestContent"
$ws.Range("J1").Value = "ResponseContent"

# Copy the header style from the existing ResponseByFields header (now K1)
# onto the two new headers so they match the banner look.
$ws.Range("K1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Re-set the values (PasteSpecial(xlPasteFormats) only touches formatting,
# but re-assert them just in case).
$ws.Range("I1").Value = "RequestContent"
$ws.Range("J1").Value = "ResponseContent"

# --- 2. Append two new trailing columns M/N ---
$ws.Range("M1").Value = "StoreResponseVariables"
$ws.Range("N1").Value = "EvaluateFunctionVariables"

# --- 3. Add the new example row (row 6), copying row 5's look first ---
# (only A:J - K6/L6 are intentionally left blank/unset for the new row)
$ws.Range("A5:J5").Copy()
$ws.Range("A6:J6").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("A6").Value = "RiskFactor-DirtectResponse"
$ws.Range("B6").Value = "REST"
$ws.Range("D6").Value = "RiskFactor-DirtectResponse"
$ws.Range("E6").Value = "https://live.virtualandemo.com/api/riskfactor/compute"
$ws.Range("F6").Value = "POST"
$ws.Range("G6").Value = "application/json"
$ws.Range("H6").Value = 200
$ws.Range("I6").Value = "{`n  ""birthday"" : ""1918-10-24"",`n  ""postalCode"" : ""60563""`n}"
$ws.Range("J6").Value = 65
$ws.Range("M6").Value = "riskFactor=."
$ws.Range("N6").Value = "success=c~[riskFactor]=65"

$ws.Range("I6:J6").WrapText = $true

# --- 4. Hyperlinks ---
# New hyperlink for the RiskFactor example URL (row 6).
$riskUrl = "https://live.virtualandemo.com/api/riskfactor/compute"
$ws.Hyperlinks.Add($ws.Range("E6"), $riskUrl, "", "", $riskUrl)

# The pre-existing bgates URL example in row 4 gets a hyperlink too.
$bgatesUrl = "https://live.virtualandemo.com/api/persons/bgates"
$ws.Hyperlinks.Add($ws.Range("E4"), $bgatesUrl, "", "", $bgatesUrl)

$ws.Range("A1").Select()
